$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T for rows 2-17
# (column numbers: E=5 G=7 H=8 I=9 J=10 K=11 M=13 N=14 O=15 P=16 Q=17 R=18 S=19 T=20)
$colNums = @(5,7,8,9,10,11,13,14,15,16,17,18,19,20)

$data = @{
    2 = @(3, 1.951738, 5.855214, 0.1200622759026226, 0.1200622759026226, 3, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 230.3968397574353, 2073.571557816918, 0.05591559823893365, 0.05591559823893365)
    3 = @(3, 1.951738, 5.855214, 0.1200622759026226, 0.1200622759026226, 3, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 115.8666152739493, 1042.799537465544, 0.02811996516872433, 0.02811996516872434)
    4 = @(3, 1.951738, 5.855214, 0.1200622759026226, 0.1200622759026226, 3, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 110.9466721718113, 998.5200495463019, 0.0269259315954028, 0.0269259315954028)
    5 = @(3, 1.951738, 5.855214, 0.1200622759026226, 0.1200622759026226, 3, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 37.49921711691334, 337.49295405222, 0.009100780899561806, 0.009100780899561808)
    6 = @(3, 9.654910333333332, 28.964731, 0.5939273141455197, 0.5939273141455197, 3, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 1139.733319196227, 10257.59987276605, 0.2766047938973343, 0.2766047938973344)
    7 = @(3, 9.654910333333332, 28.964731, 0.5939273141455197, 0.5939273141455197, 3, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 573.1721066540749, 5158.548959886675, 0.1391046043477608, 0.1391046043477608)
    8 = @(3, 9.654910333333332, 28.964731, 0.5939273141455197, 0.5939273141455197, 3, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 548.833999030898, 4939.505991278083, 0.1331979267683884, 0.1331979267683885)
    9 = @(3, 9.654910333333332, 28.964731, 0.5939273141455197, 0.5939273141455197, 3, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 185.5021415958477, 1669.51927436263, 0.04501998913203611, 0.04501998913203612)
    10 = @(3, 1.050938666666667, 3.152816, 0.06464909129917419, 0.06464909129917419, 3, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 124.0601697455769, 1116.541527710192, 0.03010847985697565, 0.03010847985697565)
    11 = @(3, 1.050938666666667, 3.152816, 0.06464909129917419, 0.06464909129917419, 3, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 62.38988335892621, 561.5089502303359, 0.01514156034320808, 0.01514156034320809)
    12 = @(3, 1.050938666666667, 3.152816, 0.06464909129917419, 0.06464909129917419, 3, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 59.74067611705422, 537.666085053488, 0.01449861746281032, 0.01449861746281032)
    13 = @(3, 1.050938666666667, 3.152816, 0.06464909129917419, 0.06464909129917419, 3, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 20.19194033107556, 181.72746297968, 0.004900433636180139, 0.00490043363618014)
    14 = @(3, 3.59846, 10.79538, 0.2213613186526835, 0.2213613186526835, 3, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 424.7874520010066, 3823.08706800906, 0.1030927530431201, 0.1030927530431201)
    15 = @(3, 3.59846, 10.79538, 0.2213613186526835, 0.2213613186526835, 3, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 213.6256917673866, 1922.63122590648, 0.05184536544405435, 0.05184536544405436)
    16 = @(3, 3.59846, 10.79538, 0.2213613186526835, 0.2213613186526835, 3, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 204.5546902009266, 1840.99221180834, 0.04964390087644609, 0.04964390087644609)
    17 = @(3, 3.59846, 10.79538, 0.2213613186526835, 0.2213613186526835, 3, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 69.13808760526666, 622.2427884474, 0.01677929928906297, 0.01677929928906297)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $colNums.Length; $i++) {
        $ws.Cells.Item([int]$row, $colNums[$i]).Value = $vals[$i]
    }
}
